$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global_variable")

# Row 24 ("SCPT location") previously listed Taichung boreholes; replace with
# the Hsinchu boreholes used for the integrated lab-test / PSD plotting run.
# The leading "'" mirrors the original cells' quote-prefixed (force-text)
# formatting so the saved style matches the source cells exactly.
$ws.Range("B24").Value = "'BH01-HSINCHU"
$ws.Range("C24").Value = "'BH03-HSINCHU"
$ws.Range("D24").Value = "'BH04-HSINCHU"
$ws.Range("I21").Select() | Out-Null
